$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace straight/curly double quotes around certain quoted words with single quotes,
# in the English (column C) dialogue lines, per commit "update on 20210731".

$ws.Range("C42").Value = "[name=`"Strange-Looking Tourist`"]   Don't blame me. That's right, it's all your fault. Just how many singers were sacrificed for your so-called 'rap empire'...`n"

$ws.Range("C45").Value = "[name=`"Emperor`"]   You think I'd forget who you are just because you decided to roll around in the sewers? Will pointing your weapon at these youngsters satisfy you, 'Songster?'`n"

$ws.Range("C47").Value = "[name=`"Emperor`"]   If I wasn’t such a juicy target for you, how many rising stars would you have continued to 'eliminate' just for not suiting your tastes?`n"

$ws.Range("C59").Value = "[name=`"Emperor`"]   My past is more expansive than this sea of falsity; so where are these 'enemies' that await my majesty?`n"

$ws.Range("C65").Value = "[name=`"Sora`"]   Tex- err, I mean, 'Sora,' that's… because you stayed in the hotel reading books the whole time…`n"
